$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 712.75
$ws.Range("I41").Value = 511
$ws.Range("J41").Value = 1049
$ws.Range("K41").Value = 511
$ws.Range("L41").Value = 1049
$ws.Range("M41").Value = -71
$ws.Range("N41").Value = -1929

$ws.Range("H43").Value = 2332.6667
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H106").Value = 96377.63
$ws.Range("I106").Value = 129771.25
$ws.Range("J106").Value = 7328
$ws.Range("K106").Value = 129771.25
$ws.Range("L106").Value = 7328
$ws.Range("M106").Value = -129140.25
$ws.Range("N106").Value = -8590

$ws.Range("H135").Value = 1670.2106
$ws.Range("I135").Value = 972.8461
$ws.Range("K135").Value = 8755.6149
$ws.Range("M135").Value = -6220.6149

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3574763
$ws.Range("I2").Value = 7144607
$ws.Range("K2").Value = 7144607
$ws.Range("M2").Value = -7144494

$ws.Range("H45").Value = 7016.25
$ws.Range("I45").Value = 9007.941000000001
$ws.Range("K45").Value = 9007.941000000001
$ws.Range("M45").Value = -8630.941000000001

$ws.Range("H61").Value = 3156.7693
$ws.Range("I61").Value = 3187.2632
$ws.Range("K61").Value = 3187.2632
$ws.Range("M61").Value = -2975.2632

$ws.Range("H116").Value = 3574763
$ws.Range("I116").Value = 7144607
$ws.Range("K116").Value = 7144607
$ws.Range("M116").Value = -7142313

$ws.Range("H122").Value = 5844.4644
$ws.Range("J122").Value = 5488.6
$ws.Range("L122").Value = 16465.8
$ws.Range("N122").Value = -21365.8

$ws.Range("H132").Value = 1478.9143
$ws.Range("I132").Value = 1445.0333
$ws.Range("K132").Value = 4335.0999
$ws.Range("M132").Value = -1805.0999

$ws.Range("H136").Value = 3156.7693
$ws.Range("I136").Value = 3187.2632
$ws.Range("K136").Value = 9561.7896
$ws.Range("M136").Value = -7011.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3574763
$ws.Range("I3").Value = 7144607
$ws.Range("K3").Value = 7144607
$ws.Range("M3").Value = -7144493

$ws.Range("H20").Value = 1758.7018
$ws.Range("I20").Value = 1715.1522
$ws.Range("J20").Value = 1940.8182
$ws.Range("K20").Value = 1715.1522
$ws.Range("L20").Value = 1940.8182
$ws.Range("M20").Value = -1468.1522
$ws.Range("N20").Value = -2434.8182

$ws.Range("H97").Value = 7133.3335
$ws.Range("I97").Value = 7133.3335
$ws.Range("K97").Value = 7133.3335
$ws.Range("M97").Value = -6142.3335

$ws.Range("H134").Value = 2781.6226
$ws.Range("I134").Value = 2483.6365
$ws.Range("K134").Value = 7450.9095
$ws.Range("M134").Value = -4915.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2484.2632
$ws.Range("I31").Value = 1796.2
$ws.Range("J31").Value = 3248.7778
$ws.Range("K31").Value = 1796.2
$ws.Range("L31").Value = 3248.7778
$ws.Range("M31").Value = -1501.2
$ws.Range("N31").Value = -3838.7778

$ws.Range("H34").Value = 2484.2632
$ws.Range("I34").Value = 1796.2
$ws.Range("J34").Value = 3248.7778
$ws.Range("K34").Value = 1796.2
$ws.Range("L34").Value = 3248.7778
$ws.Range("M34").Value = -1594.2
$ws.Range("N34").Value = -3652.7778

$ws.Range("H43").Value = 241666.33
$ws.Range("J43").Value = 241666.33
$ws.Range("L43").Value = 241666.33
$ws.Range("N43").Value = -242034.33

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H101").Value = 241666.33
$ws.Range("J101").Value = 241666.33
$ws.Range("L101").Value = 241666.33
$ws.Range("N101").Value = -248156.33

$ws.Range("H105").Value = 2562.24
$ws.Range("I105").Value = 1236.909
$ws.Range("J105").Value = 12281.333
$ws.Range("K105").Value = 1236.909
$ws.Range("L105").Value = 12281.333
$ws.Range("M105").Value = 510.0909999999999
$ws.Range("N105").Value = -15775.333

$ws.Range("H122").Value = 12718.723
$ws.Range("I122").Value = 14525.923
$ws.Range("J122").Value = 8020
$ws.Range("K122").Value = 43577.769
$ws.Range("L122").Value = 24060
$ws.Range("M122").Value = -41127.769
$ws.Range("N122").Value = -28960

$ws.Range("H134").Value = 2652.7878
$ws.Range("I134").Value = 2329.0344
$ws.Range("K134").Value = 6987.1032
$ws.Range("M134").Value = -4452.1032

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 425
$ws.Range("J32").Value = 425
$ws.Range("L32").Value = 1275
$ws.Range("N32").Value = -1841

$ws.Range("H46").Value = 535.2857
$ws.Range("J46").Value = 799.5
$ws.Range("L46").Value = 2398.5
$ws.Range("N46").Value = -2580.5

$ws.Range("H129").Value = 1130.2
$ws.Range("I129").Value = 783.1
$ws.Range("J129").Value = 1824.4
$ws.Range("K129").Value = 2349.3
$ws.Range("L129").Value = 5473.200000000001
$ws.Range("M129").Value = 2650.7
$ws.Range("N129").Value = -15473.2

$ws.Range("H134").Value = 2205.4707
$ws.Range("I134").Value = 1037.9231
$ws.Range("K134").Value = 3113.7693
$ws.Range("M134").Value = 1956.2307

$ws.Range("H139").Value = 5425.875
$ws.Range("I139").Value = 1562
$ws.Range("K139").Value = 4686
$ws.Range("M139").Value = 454

$ws.Range("H140").Value = 1913.1111
$ws.Range("I140").Value = 1913.1111
$ws.Range("K140").Value = 5739.3333
$ws.Range("M140").Value = -559.3333000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7212.6665
$ws.Range("I70").Value = 6864.25
$ws.Range("K70").Value = 6864.25
$ws.Range("M70").Value = -6594.25

$ws.Range("H73").Value = 7212.6665
$ws.Range("I73").Value = 6864.25
$ws.Range("K73").Value = 6864.25
$ws.Range("M73").Value = -5928.25

$ws.Range("H95").Value = 29333.334
$ws.Range("J95").Value = 29333.334
$ws.Range("L95").Value = 29333.334
$ws.Range("N95").Value = -34825.334

$ws.Range("H99").Value = 13513.6
$ws.Range("I99").Value = 7190
$ws.Range("J99").Value = 22999
$ws.Range("K99").Value = 7190
$ws.Range("L99").Value = 22999
$ws.Range("M99").Value = -4944
$ws.Range("N99").Value = -27491

$ws.Range("H102").Value = 8784.909
$ws.Range("I102").Value = 10803.667
$ws.Range("K102").Value = 10803.667
$ws.Range("M102").Value = -9181.666999999999

$ws.Range("H126").Value = 4779.778
$ws.Range("I126").Value = 4127.5
$ws.Range("K126").Value = 12382.5
$ws.Range("M126").Value = -9912.5

$ws.Range("H136").Value = 81333
$ws.Range("J136").Value = 81333
$ws.Range("L136").Value = 243999
$ws.Range("N136").Value = -249099

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3250.158
$ws.Range("I46").Value = 2749.8333
$ws.Range("J46").Value = 3481.077
$ws.Range("K46").Value = 2749.8333
$ws.Range("L46").Value = 3481.077
$ws.Range("M46").Value = -2561.8333
$ws.Range("N46").Value = -3857.077

$ws.Range("H55").Value = 407.10526
$ws.Range("I55").Value = 171.92857
$ws.Range("K55").Value = 171.92857
$ws.Range("M55").Value = 1.071429999999992

$ws.Range("H68").Value = 5353.4546
$ws.Range("I68").Value = 894.7368
$ws.Range("J68").Value = 33592
$ws.Range("K68").Value = 894.7368
$ws.Range("L68").Value = 33592
$ws.Range("M68").Value = -145.7368
$ws.Range("N68").Value = -35090

$ws.Range("H71").Value = 5353.4546
$ws.Range("I71").Value = 894.7368
$ws.Range("J71").Value = 33592
$ws.Range("K71").Value = 4473.684
$ws.Range("L71").Value = 167960
$ws.Range("M71").Value = -729.6840000000002
$ws.Range("N71").Value = -175448

$ws.Range("H105").Value = 83490
$ws.Range("J105").Value = 83490
$ws.Range("L105").Value = 83490
$ws.Range("N105").Value = -90478

$ws.Range("H112").Value = 54437.5
$ws.Range("J112").Value = 54437.5
$ws.Range("L112").Value = 54437.5
$ws.Range("N112").Value = -57391.5

$ws.Range("H132").Value = 49571.645
$ws.Range("I132").Value = 57091.5
$ws.Range("K132").Value = 171274.5
$ws.Range("M132").Value = -168744.5

$ws.Range("H136").Value = 2573.742
$ws.Range("I136").Value = 1710.8518
$ws.Range("K136").Value = 5132.555399999999
$ws.Range("M136").Value = -2582.555399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 7500
$ws.Range("J63").Value = 7500
$ws.Range("L63").Value = 7500
$ws.Range("N63").Value = -8748

$ws.Range("H66").Value = 7500
$ws.Range("J66").Value = 7500
$ws.Range("L66").Value = 22500
$ws.Range("N66").Value = -28740

$ws.Range("H81").Value = 2464.1667
$ws.Range("I81").Value = 2464.1667
$ws.Range("K81").Value = 4928.3334
$ws.Range("M81").Value = -3867.3334

$ws.Range("H84").Value = 2464.1667
$ws.Range("I84").Value = 2464.1667
$ws.Range("K84").Value = 24641.667
$ws.Range("M84").Value = -19337.667

$ws.Range("H122").Value = 4113.4614
$ws.Range("I122").Value = 4639.727
$ws.Range("J122").Value = 1219
$ws.Range("K122").Value = 13919.181
$ws.Range("L122").Value = 3657
$ws.Range("M122").Value = -11469.181
$ws.Range("N122").Value = -8557

$ws.Range("H132").Value = 2779.0625
$ws.Range("J132").Value = 4975.3335
$ws.Range("L132").Value = 14926.0005
$ws.Range("N132").Value = -19986.0005

$ws.Range("H136").Value = 1821.7878
$ws.Range("I136").Value = 1852.2258
$ws.Range("K136").Value = 5556.6774
$ws.Range("M136").Value = -3006.6774


Write-Host "Applied all updates"
